# Update countries & provincias Spain
#
# Refreshes the COVID "Pais" dashboard table with newer source numbers and
# bumps the "last refreshed" timestamp banner. A few countries changed
# rank in the underlying feed (their table row swapped places with the
# neighbouring row), so for those rows both the country name (column A)
# and the stats (columns B:H) are rewritten; for the rest only the stats
# (B:H) are refreshed and the country name in column A is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp banner (A1) ---------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 10 de Septiembre de 2020 a las 05:52"

# --- Helper: write the 7 stat columns (B:H) for one row ------------------
function Set-PaisStats($r, $casosTotales, $nuevosCasos, $casosActivos, $recuperados, $casosCriticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($r, 2).Value = $casosTotales
    $ws.Cells.Item($r, 3).Value = $nuevosCasos
    $ws.Cells.Item($r, 4).Value = $casosActivos
    $ws.Cells.Item($r, 5).Value = $recuperados
    $ws.Cells.Item($r, 6).Value = $casosCriticos
    $ws.Cells.Item($r, 7).Value = $muertesHoy
    $ws.Cells.Item($r, 8).Value = $muertes
}

# --- Rows whose ranking/country stayed the same: just refresh stats -----
Set-PaisStats 20  299855 196 287950 5540  0 6  6365
Set-PaisStats 33  106584 86  100248 4702  0 0  1634
Set-PaisStats 39  89691  550 18635  61139 0 5  9917
Set-PaisStats 50  65597  379 15064  48489 0 10 2044
Set-PaisStats 160 1365   4   326    1021  0 2  18
Set-PaisStats 172 628    14  266    357   0 0  5
Set-PaisStats 173 530    3   333    178   0 0  19
Set-PaisStats 184 310    0   298    12    0 0  0
Set-PaisStats 187 234    0   155    79    0 0  0

# --- Rows that swapped rank with their neighbour: name + data move ------
# Bahamas overtakes Jordania (rows 137/138)
$ws.Cells.Item(137, 1).Value = "Bahamas"
Set-PaisStats 137 2721 64 1183 1475 0 0 63

$ws.Cells.Item(138, 1).Value = "Jordania"
Set-PaisStats 138 2659 0 1919 721 0 0 19

# Birmania overtakes Yemen (rows 150/151)
$ws.Cells.Item(150, 1).Value = "Birmania"
Set-PaisStats 150 2009 120 553 1442 0 2 14

$ws.Cells.Item(151, 1).Value = "Yemen"
Set-PaisStats 151 1999 0 1209 214 0 0 576

# Montserrat overtakes Islas Malvinas (rows 214/215)
$ws.Cells.Item(214, 1).Value = "Montserrat"
Set-PaisStats 214 13 0 12 0 0 0 1

$ws.Cells.Item(215, 1).Value = "Islas Malvinas"
Set-PaisStats 215 13 0 13 0 0 0 0
